# Update the cryptos list worksheet with refreshed price / volume figures.
# All values in columns D (Price) and E (Volume(1h)) are stored as plain
# text (inlineStr) in the original workbook, so we force text assignment
# via NumberFormat "@" before writing to avoid Excel auto-converting
# values such as "1.636.40" or "26.716.61" into dates/numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    'D2'  = '26.716.61'
    'E2'  = '  -0.05%  '
    'D3'  = '1.635.87'
    'E3'  = '  -0.66%  '
    'E4'  = '  +0.19%  '
    'D5'  = '217.11'
    'E5'  = '  +0.53%  '
    'E7'  = '  +0.22%  '
    'E8'  = '  -0.65%  '
    'E9'  = '  -0.79%  '
    'D10' = '19.02'
    'E10' = '  -0.76%  '
    'D11' = '0.0844'
    'E11' = '  +0.22%  '
    'D12' = '1.863.77'
    'E12' = '  -0.65%  '
    'D13' = '1.642.82'
    'E14' = '  -1.14%  '
    'E15' = '  -1.39%  '
    'D16' = '64.38'
    'E16' = '  -1.38%  '
    'D17' = '26.708.94'
    'E17' = '  -0.07%  '
    'E18' = '  -2.29%  '
    'E19' = '  +0.17%  '
    'D20' = '210.51'
    'E20' = '  -3.67%  '
    'E21' = '  -0.86%  '
    'E22' = '  -1.55%  '
    'D23' = '2.33'
    'E23' = '  +1.99%  '
    'E24' = '  -2.98%  '
    'D25' = '145.68'
    'E25' = '  -0.21%  '
    'E26' = '  -0.04%  '
    'E27' = '  -2.22%  '
    'E28' = '  -1.02%  '
    'D29' = '15.55'
    'E29' = '  -1.18%  '
    'D30' = '0.0504'
    'E30' = '  -2.37%  '
    'E31' = '  +0.63%  '
    'E32' = '  -0.47%  '
    'E33' = '  -1.61%  '
    'D34' = '1.273.67'
    'E34' = '  -0.53%  '
    'E35' = '  -1.58%  '
    'D36' = '2.44'
    'E36' = '  +0.49%  '
    'E37' = '  -1.87%  '
    'E38' = '  -1.13%  '
    'E39' = '  -2.17%  '
    'E40' = '  +0.16%  '
    'E41' = '  -1.67%  '
    'E42' = '  -2.31%  '
    'D43' = '1.773.46'
    'E43' = '  -0.66%  '
    'E44' = '  -3.58%  '
    'D45' = '60.52'
    'E45' = '  +1.19%  '
    'D46' = '91.00'
    'E46' = '  -1.07%  '
    'E47' = '  -2.29%  '
    'E48' = '  +0.75%  '
    'E49' = '  -3.25%  '
    'D50' = '0.0959'
    'E50' = '  -0.81%  '
    'E51' = '  -0.09%  '
}

foreach ($key in $updates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$key]
}
